$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values on row 13
$ws.Range("C13").Value = 0.286
$ws.Range("E13").Value = 0.25

# Update the active selection to E14 (was C14)
$ws.Range("E14").Select()
